$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Add a new "Bottle" sheet at the end (after "SCOTCH") with bottle
#    prices for Whisky / Rum / Scotch / Vodka (x2 rows).
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$bottle = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$bottle.Name = "Bottle"

$bottle.Range("A1").Value = "Whisky"
$bottle.Range("B1").Value = 1000
$bottle.Range("A2").Value = "Rum"
$bottle.Range("B2").Value = 1500
$bottle.Range("A3").Value = "Scotch"
$bottle.Range("B3").Value = 1200
$bottle.Range("A4").Value = "Vodka"
$bottle.Range("B4").Value = 1500
$bottle.Range("A5").Value = "Vodka"
$bottle.Range("B5").Value = 1600

$bottle.Range("B6").Select()

# ------------------------------------------------------------------
# 2. SCOTCH sheet: widen column A (cosmetic column-width tweak).
# ------------------------------------------------------------------
$wsScotch = $wb.Worksheets.Item("SCOTCH")
$wsScotch.Columns.Item(1).ColumnWidth = 22.3

# ------------------------------------------------------------------
# 3. "FOR THE TEETOTALLERS" sheet: append a new "CL" row priced at 20,
#    matching the bold style already used for the rest of the column,
#    then leave the selection on the next empty row (A15) as the
#    final/active sheet.
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("FOR THE TEETOTALLERS")
$ws1.Range("A14").Value = "CL"
$ws1.Range("B14").Value = 20
$ws1.Range("A14:B14").Font.Bold = $true

$ws1.Activate()
$ws1.Range("A15").Select()
